$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 0.8673686043034796
$ws.Range("F1").Value = -1.570796292848413

$ws.Range("E2").Value = 0.869309796466788
$ws.Range("F2").Value = -1.570796289207822

$ws.Range("E3").Value = 0.8780076243366395
$ws.Range("F3").Value = -1.570796272895559

$ws.Range("E4").Value = 0.8902446423384895
$ws.Range("F4").Value = -1.570796249945755

$ws.Range("E5").Value = 0.8989424702083411
$ws.Range("F5").Value = -1.570796233633492

$ws.Range("E6").Value = 0.9008836623716494
$ws.Range("F6").Value = -1.570796229992901
